# Resize/reposition the icon picture and its caption textbox on the
# "Challenges and Takeaways" slide (slide 1).
#
# Target values come from the canonical OOXML (EMUs); the PowerPoint
# object model works in points, so divide by 12700 (EMU per point).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$emuPerPt = 12700

$picture = $s.Shapes.Item("Picture 6")
$picture.Left   = 2624475 / $emuPerPt
$picture.Top    = 1240117 / $emuPerPt
$picture.Width  = 3905394 / $emuPerPt
$picture.Height = 3905394 / $emuPerPt

$caption = $s.Shapes.Item("TextBox 7")
$caption.Left   = 2624475 / $emuPerPt
# tiny epsilon nudges compensate for the single-precision point storage
# used internally so the round-tripped EMU lands exactly on target
$caption.Top    = (5914014 / $emuPerPt) + 0.00001
$caption.Width  = 3905394 / $emuPerPt
$caption.Height = (461665 / $emuPerPt) + 0.00001
